$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.579.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").Value = "'1.914.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.71%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'315.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5072"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").Value = "'0.3956"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").Value = "'0.09838"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "'1.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.67%  "
$ws.Range("D11").Value = "'42.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.69%  "
$ws.Range("D12").Value = "'6.548"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("D13").Value = "'21.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.14%  "
$ws.Range("D14").Value = "'1.920.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.14%  "
$ws.Range("D15").Value = "'7.586"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.43%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "'0.00001141"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "'94.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").Value = "'0.06666"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").Value = "'18.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.86%  "
$ws.Range("D21").Value = "'0.9997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "'6.326"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.19%  "
$ws.Range("D23").Value = "'28.635.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'11.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("D25").Value = "'2.288"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D26").Value = "'2.736"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +14.78%  "
$ws.Range("D27").Value = "'2.136.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.77%  "
$ws.Range("D28").Value = "'21.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.02%  "
$ws.Range("D29").Value = "'159.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").Value = "'128.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("E31").Value = "  +7.17%  "
$ws.Range("D32").Value = "'0.1077"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").Value = "'5.754"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.80%  "
$ws.Range("D34").Value = "'3.646"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "'9.877"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.74%  "
$ws.Range("D36").Value = "'0.06813"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("D37").Value = "'0.02450"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.40%  "
$ws.Range("D38").Value = "'1.274"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.06%  "
$ws.Range("D39").Value = "'0.2235"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.71%  "
$ws.Range("D40").Value = "'11.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.95%  "
$ws.Range("D41").Value = "'5.109"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.81%  "
$ws.Range("D42").Value = "'0.6450"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.74%  "
$ws.Range("D43").Value = "'1.192"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'13.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.02%  "
$ws.Range("D46").Value = "'0.6099"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.78%  "
$ws.Range("D47").Value = "'3.812"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.32%  "
$ws.Range("D48").Value = "'1.285"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("E49").Value = "  +5.73%  "
$ws.Range("D50").Value = "'125.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("E51").Value = "  +3.40%  "
